$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily TGP (terminal gate pricing) refresh: the "latest date" rows shift
# down to become the "previous date" rows (with their old values), and a
# new set of rows for the new latest date is written in with fresh prices.

# --- New South Wales (rows 8-13) ---
$ws.Range("A8").Value = 45997
$ws.Range("D8").Value = 162.59
$ws.Range("E8").Value = 161.45
$ws.Range("F8").Value = 171.45
$ws.Range("G8").Value = 161.57

$ws.Range("A9").Value = 45997
$ws.Range("D9").Value = 162.59
$ws.Range("E9").Value = 161.45
$ws.Range("F9").Value = 171.45
$ws.Range("G9").Value = 161.57

$ws.Range("A10").Value = 45997
$ws.Range("D10").Value = 164.79
$ws.Range("E10").Value = 163.95
$ws.Range("F10").Value = 173.95
$ws.Range("G10").Value = 164.47

$ws.Range("A11").Value = 45996
$ws.Range("D11").Value = 163.05
$ws.Range("E11").Value = 161.81
$ws.Range("F11").Value = 171.81
$ws.Range("G11").Value = 161.93

$ws.Range("A12").Value = 45996
$ws.Range("D12").Value = 163.05
$ws.Range("E12").Value = 161.81
$ws.Range("F12").Value = 171.81
$ws.Range("G12").Value = 161.93

$ws.Range("A13").Value = 45996
$ws.Range("D13").Value = 165.14
$ws.Range("E13").Value = 164.43
$ws.Range("F13").Value = 174.43
$ws.Range("G13").Value = 164.95

# --- Northern Territory (rows 17-18) ---
$ws.Range("A17").Value = 45997
$ws.Range("D17").Value = 167.71
$ws.Range("E17").Value = 167.02
$ws.Range("F17").Value = 177.02

$ws.Range("A18").Value = 45996
$ws.Range("D18").Value = 168.08
$ws.Range("E18").Value = 167.19
$ws.Range("F18").Value = 177.19

# --- Queensland (rows 22-31) ---
$ws.Range("A22").Value = 45997
$ws.Range("D22").Value = 163.84
$ws.Range("E22").Value = 163.24
$ws.Range("F22").Value = 172.84
$ws.Range("G22").Value = 164.4

$ws.Range("A23").Value = 45997
$ws.Range("D23").Value = 169.82
$ws.Range("E23").Value = 168.17
$ws.Range("F23").Value = 178.17

$ws.Range("A24").Value = 45997
$ws.Range("D24").Value = 169.58
$ws.Range("E24").Value = 168.49
$ws.Range("F24").Value = 178.49

$ws.Range("A25").Value = 45997
$ws.Range("D25").Value = 170.19
$ws.Range("E25").Value = 167.93
$ws.Range("F25").Value = 177.93
$ws.Range("G25").Value = 167.7

$ws.Range("A26").Value = 45997
$ws.Range("D26").Value = 168.99
$ws.Range("E26").Value = 169.39
$ws.Range("F26").Value = 179.39

$ws.Range("A27").Value = 45996
$ws.Range("D27").Value = 164.18
$ws.Range("E27").Value = 163.72
$ws.Range("F27").Value = 173.32
$ws.Range("G27").Value = 164.88

$ws.Range("A28").Value = 45996
$ws.Range("D28").Value = 170.17
$ws.Range("E28").Value = 168.65
$ws.Range("F28").Value = 178.65

$ws.Range("A29").Value = 45996
$ws.Range("D29").Value = 169.93
$ws.Range("E29").Value = 168.97
$ws.Range("F29").Value = 178.97

$ws.Range("A30").Value = 45996
$ws.Range("D30").Value = 170.54
$ws.Range("E30").Value = 168.41
$ws.Range("F30").Value = 178.41
$ws.Range("G30").Value = 168.18

$ws.Range("A31").Value = 45996
$ws.Range("D31").Value = 169.34
$ws.Range("E31").Value = 169.87
$ws.Range("F31").Value = 179.87

# --- South Australia (rows 35-36) ---
$ws.Range("A35").Value = 45997
$ws.Range("D35").Value = 163.06
$ws.Range("E35").Value = 160.95
$ws.Range("F35").Value = 169.95

$ws.Range("A36").Value = 45996
$ws.Range("D36").Value = 163.52
$ws.Range("E36").Value = 161.43
$ws.Range("F36").Value = 170.43

# --- Tasmania (rows 40-43) ---
$ws.Range("A40").Value = 45997
$ws.Range("D40").Value = 169.09
$ws.Range("E40").Value = 167.68
$ws.Range("F40").Value = 177.68

$ws.Range("A41").Value = 45997
$ws.Range("D41").Value = 168.8
$ws.Range("E41").Value = 168.09
$ws.Range("F41").Value = 178.09

$ws.Range("A42").Value = 45996
$ws.Range("D42").Value = 169.43
$ws.Range("E42").Value = 167.81
$ws.Range("F42").Value = 177.81

$ws.Range("A43").Value = 45996
$ws.Range("D43").Value = 169.14
$ws.Range("E43").Value = 168.23
$ws.Range("F43").Value = 178.23

# --- Victoria (rows 47-50) ---
$ws.Range("A47").Value = 45997
$ws.Range("D47").Value = 163.65
$ws.Range("E47").Value = 162.9
$ws.Range("F47").Value = 172.9

$ws.Range("A48").Value = 45997
$ws.Range("D48").Value = 163.48
$ws.Range("E48").Value = 162.99
$ws.Range("F48").Value = 172.99

$ws.Range("A49").Value = 45996
$ws.Range("D49").Value = 163.68
$ws.Range("E49").Value = 163.33
$ws.Range("F49").Value = 173.33

$ws.Range("A50").Value = 45996
$ws.Range("D50").Value = 163.52
$ws.Range("E50").Value = 163.44
$ws.Range("F50").Value = 173.44

# --- Western Australia (rows 54-65) ---
$ws.Range("A54").Value = 45997
$ws.Range("D54").Value = 178.64
$ws.Range("E54").Value = 178.77
$ws.Range("F54").Value = 188.77

$ws.Range("A55").Value = 45997
$ws.Range("D55").Value = 166.76
$ws.Range("E55").Value = 173.69
$ws.Range("F55").Value = 183.69

$ws.Range("A56").Value = 45997
$ws.Range("D56").Value = 169.24

$ws.Range("A57").Value = 45997
$ws.Range("D57").Value = 168.2
$ws.Range("E57").Value = 167.96

$ws.Range("A58").Value = 45997
$ws.Range("D58").Value = 164.1
$ws.Range("E58").Value = 164.01
$ws.Range("F58").Value = 174.01

$ws.Range("A59").Value = 45997
$ws.Range("D59").Value = 170.61
$ws.Range("E59").Value = 176.02

$ws.Range("A60").Value = 45996
$ws.Range("D60").Value = 179.01
$ws.Range("E60").Value = 179.19
$ws.Range("F60").Value = 189.19

$ws.Range("A61").Value = 45996
$ws.Range("D61").Value = 167.12
$ws.Range("E61").Value = 174.3
$ws.Range("F61").Value = 184.3

$ws.Range("A62").Value = 45996
$ws.Range("D62").Value = 169.59

$ws.Range("A63").Value = 45996
$ws.Range("D63").Value = 168.57
$ws.Range("E63").Value = 168.56

$ws.Range("A64").Value = 45996
$ws.Range("D64").Value = 164.47
$ws.Range("E64").Value = 164.62
$ws.Range("F64").Value = 174.62

$ws.Range("A65").Value = 45996
$ws.Range("D65").Value = 170.98
$ws.Range("E65").Value = 176.48
